$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue/week date roll-forward) ---
$ws.Range("A8").Value = "Volume 31   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("C16").Value2 = 2
$ws.Range("E16").Value2 = 0
$ws.Range("F16").Value2 = 8
$ws.Range("G16").Value2 = 8
$ws.Range("I16").Value2 = 14
$ws.Range("J16").Value2 = 19
$ws.Range("K16").Value2 = -26.315789473684
$ws.Range("L16").Value2 = 27.272727272727
$ws.Range("M16").Value2 = 16.666666666666
$ws.Range("N16").Value2 = -87.155963302752
$ws.Range("C17").Value2 = 2
$ws.Range("D17").Value2 = 1
$ws.Range("E17").Value2 = 100
$ws.Range("F17").Value2 = 5
$ws.Range("H17").Value2 = -44.444444444444
$ws.Range("I17").Value2 = 19
$ws.Range("J17").Value2 = 19
$ws.Range("L17").Value2 = -20.833333333333
$ws.Range("M17").Value2 = 90
$ws.Range("N17").Value2 = -42.424242424242
$ws.Range("C18").Value2 = 2
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 7
$ws.Range("G18").Value2 = 12
$ws.Range("H18").Value2 = -41.666666666666
$ws.Range("I18").Value2 = 19
$ws.Range("J18").Value2 = 30
$ws.Range("K18").Value2 = -36.666666666666
$ws.Range("L18").Value2 = -53.658536585365
$ws.Range("M18").Value2 = -34.482758620689
$ws.Range("N18").Value2 = -91.479820627802
$ws.Range("C19").Value2 = 10
$ws.Range("D19").Value2 = 7
$ws.Range("E19").Value2 = 42.857142857142
$ws.Range("F19").Value2 = 40
$ws.Range("G19").Value2 = 50
$ws.Range("H19").Value2 = -20
$ws.Range("I19").Value2 = 118
$ws.Range("J19").Value2 = 114
$ws.Range("K19").Value2 = 3.508771929824
$ws.Range("L19").Value2 = -7.8125
$ws.Range("M19").Value2 = -4.065040650406
$ws.Range("N19").Value2 = -65.895953757225
$ws.Range("G20").Value2 = 2
$ws.Range("H20").Value2 = 50
$ws.Range("L20").Value2 = 20
$ws.Range("N20").Value2 = -95.121951219512
$ws.Range("C21").Value2 = 16
$ws.Range("D21").Value2 = 12
$ws.Range("E21").Value2 = 33.333333333333
$ws.Range("F21").Value2 = 65
$ws.Range("G21").Value2 = 81
$ws.Range("H21").Value2 = -19.753086419753
$ws.Range("I21").Value2 = 180
$ws.Range("J21").Value2 = 189
$ws.Range("K21").Value2 = -4.761904761904
$ws.Range("L21").Value2 = -14.691943127962
$ws.Range("M21").Value2 = 2.272727272727
$ws.Range("N21").Value2 = -78.494623655914
$ws.Range("F22").Value2 = 6
$ws.Range("H22").Value2 = 200
$ws.Range("I22").Value2 = 10
$ws.Range("J22").Value2 = 6
$ws.Range("K22").Value2 = 66.666666666666
$ws.Range("L22").Value2 = 66.666666666666
$ws.Range("M22").Value2 = 25
$ws.Range("C24").Value2 = 24
$ws.Range("D24").Value2 = 13
$ws.Range("E24").Value2 = 84.615384615384
$ws.Range("G24").Value2 = 63
$ws.Range("H24").Value2 = 19.047619047619
$ws.Range("I24").Value2 = 188
$ws.Range("J24").Value2 = 172
$ws.Range("K24").Value2 = 9.302325581395
$ws.Range("L24").Value2 = -8.737864077669
$ws.Range("M24").Value2 = 57.983193277310
$ws.Range("C25").Value2 = 24
$ws.Range("D25").Value2 = 11
$ws.Range("E25").Value2 = 118.181818181818
$ws.Range("F25").Value2 = 69
$ws.Range("G25").Value2 = 49
$ws.Range("H25").Value2 = 40.816326530612
$ws.Range("I25").Value2 = 156
$ws.Range("J25").Value2 = 135
$ws.Range("K25").Value2 = 15.555555555555
$ws.Range("L25").Value2 = -8.235294117647
$ws.Range("C26").Value2 = 3
$ws.Range("D26").Value2 = 9
$ws.Range("E26").Value2 = -66.666666666666
$ws.Range("F26").Value2 = 13
$ws.Range("G26").Value2 = 21
$ws.Range("H26").Value2 = -38.095238095238
$ws.Range("I26").Value2 = 38
$ws.Range("J26").Value2 = 48
$ws.Range("K26").Value2 = -20.833333333333
$ws.Range("L26").Value2 = -2.564102564102
$ws.Range("M26").Value2 = -17.391304347826
$ws.Range("F28").Value2 = 5
$ws.Range("G28").Value2 = 2
$ws.Range("H28").Value2 = 150
$ws.Range("I28").Value2 = 10
$ws.Range("K28").Value2 = 25
$ws.Range("L28").Value2 = -16.666666666666
$ws.Range("H31").Value2 = -100

# --- Cells changing from a number to the "0" placeholder text ---
# (reuses the shared string already used by neighboring dash cells)
$ws.Range("C15").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C20").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("F31").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("F31").PasteSpecial(-4122)

# --- Cells changing from a number to the "***.*" placeholder text ---
$ws.Range("E28").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# --- Cells changing from placeholder text to an actual number ---
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value2 = 1
$ws.Range("N15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value2 = 0
$ws.Range("N15").Copy()
$ws.Range("L29").PasteSpecial(-4122)
$ws.Range("L29").Value2 = -100
$ws.Range("N15").Copy()
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value2 = -100

$excel.CutCopyMode = 0
